# Add arctic fox rabies codon-usage sequences (rows 428-432) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @("AY352487",10,17,4,9,3,3,5,11,9,9,9,6,9,6,6,7,5,10,4,4,5,7,5,5,0,8,8,4,9,12,6,15,1,15,6,6,7,6,4,10,10,7,19,12,11,9,23,4,2,3,0,1,3,12,4,5,4,8,14),
    @("EF611828",10,17,4,10,3,3,5,10,10,8,9,6,9,6,6,6,5,10,5,5,4,7,5,5,0,9,7,5,9,11,7,14,1,15,6,6,7,6,4,9,11,7,19,12,11,9,23,4,2,3,0,1,3,12,4,5,4,8,14),
    @("EF611830",10,17,4,9,3,3,5,11,9,9,9,6,9,6,6,7,5,10,4,4,5,7,5,5,0,8,8,5,9,12,6,14,1,15,6,6,7,6,4,10,10,6,20,12,11,9,23,4,2,3,0,1,3,12,4,5,4,7,15),
    @("EF611851",10,17,2,10,4,3,6,10,10,8,9,8,7,6,6,8,4,10,4,5,4,7,5,5,0,12,5,5,8,11,7,14,1,16,5,6,7,6,4,12,8,5,21,12,11,11,21,4,2,3,0,1,3,13,3,7,3,9,12),
    @("L20675",10,17,4,9,3,3,4,12,11,7,9,7,8,6,6,7,6,9,3,4,5,7,5,5,1,10,6,6,7,11,7,14,2,14,7,6,7,6,4,9,12,7,18,12,11,11,21,4,2,3,0,1,3,12,4,5,3,10,13)
)

$startRow = 428
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value = $rowVals[$c]
    }
}

$ws.Range("F433").Select() | Out-Null
